# Commit: Remove quiz module and change IP logic for multiple people to use
# the same IP (bayreuth). Observable workbook effect: three new worksheets
# (Round7, Round8, Round9) are appended, each holding a single column of
# usernames (the new per-person-same-IP logic no longer needs a separate
# IP column per round), and a couple of leftover UI-selection tweaks on the
# pre-existing Round2/Round3 sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the three new worksheets at the end of the tab strip.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$round7 = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$round7.Name = "Round7"

$round8 = $wb.Worksheets.Add([Type]::Missing, $round7)
$round8.Name = "Round8"

$round9 = $wb.Worksheets.Add([Type]::Missing, $round8)
$round9.Name = "Round9"

# ---------------------------------------------------------------------
# 2. Round7 data - fresh "spring" themed usernames.
# ---------------------------------------------------------------------
$round7Names = @(
    "LenzLena",
    "BasilikumBasti",
    "SonnenscheinStefan",
    "FrühlingsFabi",
    "LeberblümchenLudwig",
    "ApfelblüteAlbert",
    "DuftDaniel",
    "GrünGloriya",
    "VogelViktor",
    "MaikäferManfred",
    "FredFruchtig",
    "FlanierFlo",
    "SommerSebastian"
)
for ($i = 0; $i -lt $round7Names.Count; $i++) {
    $round7.Cells.Item($i + 1, 1).Value = $round7Names[$i]
}
$round7.Columns.Item(1).ColumnWidth = 25.25

# ---------------------------------------------------------------------
# 3. Round8 data.
# ---------------------------------------------------------------------
$round8Names = @(
    "KnospenKorbinian",
    "SpargelSabrina",
    "PusteblumenPhilip",
    "HellHerbert",
    "GezwitscherGeli",
    "AprilwetterAndreas",
    "BlütenBene",
    "EisheiligenEsther",
    "RadlRudi",
    "PaprikaPetra",
    "WeidenkätzchenWerner",
    "SonnenscheinSascha",
    "AperolAlex"
)
for ($i = 0; $i -lt $round8Names.Count; $i++) {
    $round8.Cells.Item($i + 1, 1).Value = $round8Names[$i]
}
$round8.Columns.Item(1).ColumnWidth = 28.25

# ---------------------------------------------------------------------
# 4. Round9 data - mixes freshly added names with four reused ones.
# ---------------------------------------------------------------------
$round9Names = @(
    "ErdbeerEva",
    "MarienkäferMichi",
    "RosmarinRalf",
    "BienenBernd",
    "RosenRosa",
    "WanderWolfgang",
    "TauwetterTom",
    "UnixUlla",
    "WifiWalter",
    "SudoSanta",
    "PasswordPaul"
)
for ($i = 0; $i -lt $round9Names.Count; $i++) {
    $round9.Cells.Item($i + 1, 1).Value = $round9Names[$i]
}

# ---------------------------------------------------------------------
# 4b. Match the "2 cm" page margin convention already used by the other
#     sheets in this workbook (Round2-Round4).
# ---------------------------------------------------------------------
foreach ($sheet in @($round7, $round8, $round9)) {
    $sheet.PageSetup.TopMargin = 56.692913399999995
    $sheet.PageSetup.BottomMargin = 56.692913399999995
}

# ---------------------------------------------------------------------
# 5. Leftover UI selection state on the pre-existing sheets, matching
#    where the author had last clicked before saving.
# ---------------------------------------------------------------------
$round3 = $wb.Worksheets.Item("Round3")
$round3.Range("B17:B20").Select()

$round8.Range("C29").Select()

$round9.Range("A8:A11").Select()

# Round7 ends up the active tab / last selection when the workbook is saved.
$round7.Range("A1:A13").Select()
